$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 77, pushing existing rows 77-78 down to 79-80
$ws.Range("A77:A78").EntireRow.Insert()

# Fill new row 77 (Primera, week of 2022-01-07)
$ws.Cells.Item(77, 1).Value2 = 8
$ws.Cells.Item(77, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(77, 3).Value = "Coquimbo"
$ws.Cells.Item(77, 4).Value2 = 44568
$ws.Cells.Item(77, 5).Value2 = 4
$ws.Cells.Item(77, 6).Value2 = 100112028
$ws.Cells.Item(77, 7).Value = "Sandia"
$ws.Cells.Item(77, 8).Value = "Sin especificar"
$ws.Cells.Item(77, 9).Value = "Primera"
$ws.Cells.Item(77, 10).Value2 = 2000
$ws.Cells.Item(77, 11).Value2 = 2800
$ws.Cells.Item(77, 12).Value2 = 3000
$ws.Cells.Item(77, 13).Value2 = 2900
$ws.Cells.Item(77, 14).Value = "`$/unidad"
$ws.Cells.Item(77, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(77, 16).Value2 = 2900
$ws.Cells.Item(77, 17).Value2 = 1
$ws.Cells.Item(77, 18).Value = "Hortaliza"

# Fill new row 78 (Segunda, week of 2022-01-07)
$ws.Cells.Item(78, 1).Value2 = 8
$ws.Cells.Item(78, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(78, 3).Value = "Coquimbo"
$ws.Cells.Item(78, 4).Value2 = 44568
$ws.Cells.Item(78, 5).Value2 = 4
$ws.Cells.Item(78, 6).Value2 = 100112028
$ws.Cells.Item(78, 7).Value = "Sandia"
$ws.Cells.Item(78, 8).Value = "Sin especificar"
$ws.Cells.Item(78, 9).Value = "Segunda"
$ws.Cells.Item(78, 10).Value2 = 1600
$ws.Cells.Item(78, 11).Value2 = 2000
$ws.Cells.Item(78, 12).Value2 = 2500
$ws.Cells.Item(78, 13).Value2 = 2250
$ws.Cells.Item(78, 14).Value = "`$/unidad"
$ws.Cells.Item(78, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(78, 16).Value2 = 2250
$ws.Cells.Item(78, 17).Value2 = 1
$ws.Cells.Item(78, 18).Value = "Hortaliza"
